# Daily attendance processing - 2025-11-26 21:46:42
# Reverses the order of the comma-separated "Recorded By" values in column G
# for every data row of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value2
    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
